$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3800727954645477
$ws.Range("J2").Value = 0.47906870267432
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 0.3056135938161111
$ws.Range("R2").Value = 2.750522344345
$ws.Range("S2").Value = 0.002928360305101355
$ws.Range("T2").Value = 0.004064743931158092

$ws.Range("I3").Value = 0.3800727954645477
$ws.Range("J3").Value = 0.47906870267432
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("S3").Value = 0.2709949124246149
$ws.Range("T3").Value = 0.3761575799718901

$ws.Range("I4").Value = 0.3800727954645477
$ws.Range("J4").Value = 0.47906870267432
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 0.07241783600777778
$ws.Range("R4").Value = 0.65176052407
$ws.Range("S4").Value = 0.0006939007970768366
$ws.Range("T4").Value = 0.0009631769180965627

$ws.Range("I5").Value = 0.3800727954645477
$ws.Range("J5").Value = 0.47906870267432
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 10.93865085835917
$ws.Range("R5").Value = 65.63190515015499
$ws.Range("S5").Value = 0.1048131091454519
$ws.Range("T5").Value = 0.09699135464139097

$ws.Range("I6").Value = 0.3800727954645477
$ws.Range("J6").Value = 0.47906870267432
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 0.06705480988333333
$ws.Range("R6").Value = 0.6034932889500001
$ws.Range("S6").Value = 0.0006425127923027611
$ws.Range("T6").Value = 0.0008918472117841707

$ws.Range("G7").Value = 0.755484
$ws.Range("H7").Value = 1.510968
$ws.Range("I7").Value = 0.6199272045354524
$ws.Range("J7").Value = 0.52093129732568
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 0.498478668148
$ws.Range("R7").Value = 2.990872008888
$ws.Range("S7").Value = 0.004776375050982572
$ws.Range("T7").Value = 0.004419934588785594

$ws.Range("G8").Value = 0.755484
$ws.Range("H8").Value = 1.510968
$ws.Range("I8").Value = 0.6199272045354524
$ws.Range("J8").Value = 0.52093129732568
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("Q8").Value = 46.129973413784
$ws.Range("R8").Value = 276.779840482704
$ws.Range("S8").Value = 0.4420130051596698
$ws.Range("T8").Value = 0.40902746315878

$ws.Range("G9").Value = 0.755484
$ws.Range("H9").Value = 1.510968
$ws.Range("I9").Value = 0.6199272045354524
$ws.Range("J9").Value = 0.52093129732568
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 0.118118916088
$ws.Range("R9").Value = 0.708713496528
$ws.Range("S9").Value = 0.001131804187224156
$ws.Range("T9").Value = 0.001047342476553496

$ws.Range("G10").Value = 0.755484
$ws.Range("H10").Value = 1.510968
$ws.Range("I10").Value = 0.6199272045354524
$ws.Range("J10").Value = 0.52093129732568
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 17.841759075978
$ws.Range("R10").Value = 71.367036303912
$ws.Range("S10").Value = 0.1709580336361383
$ws.Range("T10").Value = 0.1054667773550286

$ws.Range("G11").Value = 0.755484
$ws.Range("H11").Value = 1.510968
$ws.Range("I11").Value = 0.6199272045354524
$ws.Range("J11").Value = 0.52093129732568
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 0.10937141868
$ws.Range("R11").Value = 0.6562285120800001
$ws.Range("S11").Value = 0.001047986501437649
$ws.Range("T11").Value = 0.0009697797465322138
